$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (entered first by the author, hence its strings land earliest in
# the shared-string table) - "marcius 8" work-log entry.
$ws.Range("A6").Value = "március 8"
$ws.Range("B6").Value = "mindenki"

# Row 3 - "februar 14": team formation.
$ws.Range("A3").Value = "február 14"
$ws.Range("B3").Value = "mindenki"
$ws.Range("C3").Value = "csapat kialakítása"
$ws.Range("D3").Value = "megbeszéltük, hogy mindkettőnknek ötös a célja, így jól fog menni a csapatmunka"

# Row 4 - "februar 18": team name.
$ws.Range("A4").Value = "február 18"
$ws.Range("B4").Value = "mindenki"
$ws.Range("D4").Value = "-"
$ws.Range("C4").Value = "csapatnév kitalálása és regisztrálása"

# Row 5 - "marcius 5": project topic + plan.
$ws.Range("A5").Value = "március 5"
$ws.Range("B5").Value = "mindenki"
$ws.Range("C5").Value = "projekttéma kiválasztása"
$ws.Range("D5").Value = "-"

# Back to row 6 - fill remaining columns.
$ws.Range("C6").Value = "projekt tématerv megírása és beküldése"
$ws.Range("D6").Value = "összeültünk egy órára, és közösen megbeszéltük az alapvető funkciókat, oldalakat"

# Row 7 - "marcius 19": GitHub repo.
$ws.Range("A7").Value = "március 19"
$ws.Range("B7").Value = "Mindszenti Gergő"
$ws.Range("C7").Value = "GitHub repo létrehozása"
$ws.Range("D7").Value = "-"

# Row 8 - "marcius 19": register page.
$ws.Range("A8").Value = "március 19"
$ws.Range("B8").Value = "Mindszenti Gergő"
$ws.Range("C8").Value = "regisztrációs/bejelentkező oldal megalkotása"

# Leave D8 empty, matching the target state.

# Final selection, as captured in the saved workbook.
$ws.Range("C8").Select() | Out-Null
